$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,14
$data[0,0] = 12.04665772482491
$data[0,1] = 10.47893012220653
$data[0,2] = 0
$data[0,3] = 16.57036778080546
$data[0,4] = 36.40293817965991
$data[0,5] = 29.33399675307534
$data[0,6] = 14.57758611867002
$data[0,7] = 23.37023978080238
$data[0,8] = 7.616124039007515
$data[0,9] = 8.381945459088373
$data[0,10] = 12.7632868340245
$data[0,11] = 0
$data[0,12] = 18.84914455775374
$data[0,13] = 22.23839176888653
$data[1,0] = 11.79026408353352
$data[1,1] = 10.49642845510756
$data[1,2] = 0
$data[1,3] = 16.56679092471172
$data[1,4] = 36.4216525742797
$data[1,5] = 29.43497388215109
$data[1,6] = 14.62131112920121
$data[1,7] = 23.45491895119098
$data[1,8] = 7.602118139890761
$data[1,9] = 8.201369063723284
$data[1,10] = 12.74025727725266
$data[1,11] = 0
$data[1,12] = 18.89837454619737
$data[1,13] = 22.31451187347675
$data[2,0] = 11.63158686337303
$data[2,1] = 10.50781440913358
$data[2,2] = 0
$data[2,3] = 16.56717180848327
$data[2,4] = 36.44095580914389
$data[2,5] = 29.5044819586561
$data[2,6] = 14.6500107947666
$data[2,7] = 23.51056142555085
$data[2,8] = 7.593367498513179
$data[2,9] = 8.089310638353039
$data[2,10] = 12.72793029447649
$data[2,11] = 0
$data[2,12] = 18.93013295580479
$data[2,13] = 22.36504039120009
$data[3,0] = 11.56669887826899
$data[3,1] = 10.51261613069463
$data[3,2] = 0
$data[3,3] = 16.56797711741254
$data[3,4] = 36.45078759346293
$data[3,5] = 29.5346895615683
$data[3,6] = 14.66217241553606
$data[3,7] = 23.53415416242761
$data[3,8] = 7.589763910827239
$data[3,9] = 8.043408731532651
$data[3,10] = 12.7233667374978
$data[3,11] = 0
$data[3,12] = 18.94346073077387
$data[3,13] = 22.38658399963995
$data[4,0] = 11.55591329326378
$data[4,1] = 10.51342324139014
$data[4,2] = 0
$data[4,3] = 16.56815015775642
$data[4,4] = 36.45253888197637
$data[4,5] = 29.53981904707583
$data[4,6] = 14.66422001974058
$data[4,7] = 23.53812717336204
$data[4,8] = 7.589163274610927
$data[4,9] = 8.035774239913575
$data[4,10] = 12.72263683847591
$data[4,11] = 0
$data[4,12] = 18.94569713802521
$data[4,13] = 22.39021883859228
$data[5,0] = 11.63071255593971
$data[5,1] = 10.50787851090458
$data[5,2] = 0
$data[5,3] = 16.56718003440831
$data[5,4] = 36.44108044482032
$data[5,5] = 29.50488173408604
$data[5,6] = 14.65017292205599
$data[5,7] = 23.51087588808318
$data[5,8] = 7.59331905120136
$data[5,9] = 8.088692468016125
$data[5,10] = 12.72786688224492
$data[5,11] = 0
$data[5,12] = 18.9303111347706
$data[5,13] = 22.36532707829345
$data[6,0] = 11.95856502258511
$data[6,1] = 10.48483058438295
$data[6,2] = 0
$data[6,3] = 16.56860091979718
$data[6,4] = 36.40777028708906
$data[6,5] = 29.36725228498709
$data[6,6] = 14.59227834543226
$data[6,7] = 23.39868000900732
$data[6,8] = 7.611326334641016
$data[6,9] = 8.31996399077747
$data[6,10] = 12.75497228472486
$data[6,11] = 0
$data[6,12] = 18.86580189246346
$data[6,13] = 22.26385108329823
$data[7,0] = 12.58782680367744
$data[7,1] = 10.44470702974348
$data[7,2] = 0
$data[7,3] = 16.59173524191397
$data[7,4] = 36.40435909002926
$data[7,5] = 29.15716172051026
$data[7,6] = 14.49342257856
$data[7,7] = 23.20760129183498
$data[7,8] = 7.645427056532669
$data[7,9] = 8.761559265798903
$data[7,10] = 12.82233350930811
$data[7,11] = 0
$data[7,12] = 18.75140193843353
$data[7,13] = 22.09494643701044
$data[8,0] = 13.03708693595075
$data[8,1] = 10.41829356281195
$data[8,2] = 0
$data[8,3] = 16.6209887318559
$data[8,4] = 36.43943910684293
$data[8,5] = 29.03956267512915
$data[8,6] = 14.42970914719585
$data[8,7] = 23.08483275236831
$data[8,8] = 7.66972516260951
$data[8,9] = 9.075528105707001
$data[8,10] = 12.88022762845252
$data[8,11] = 0
$data[8,12] = 18.67466625804992
$data[8,13] = 21.98921394788883
$data[9,0] = 13.23771506253384
$data[9,1] = 10.40693717597245
$data[9,2] = 0
$data[9,3] = 16.63692125966396
$data[9,4] = 36.46350706440653
$data[9,5] = 28.99410240675459
$data[9,6] = 14.40265392710942
$data[9,7] = 23.03280289744376
$data[9,8] = 7.680611215831518
$data[9,9] = 9.215470601443533
$data[9,10] = 12.90832918800364
$data[9,11] = 0
$data[9,12] = 18.64133181542059
$data[9,13] = 21.94510454597169
$data[10,0] = 13.31307942776147
$data[10,1] = 10.4027311626002
$data[10,2] = 0
$data[10,3] = 16.64332832487971
$data[10,4] = 36.47378105136779
$data[10,5] = 28.97804757019926
$data[10,6] = 14.39268566600325
$data[10,7] = 23.01364935854488
$data[10,8] = 7.684709063349812
$data[10,9] = 9.268001542378244
$data[10,10] = 12.91921873161575
$data[10,11] = 0
$data[10,12] = 18.62893413405457
$data[10,13] = 21.928975442174
$data[11,0] = 13.29687653393589
$data[11,1] = 10.40363281040958
$data[11,2] = 0
$data[11,3] = 16.64193188750683
$data[11,4] = 36.47151688374908
$data[11,5] = 28.98145361010875
$data[11,6] = 14.39482019650646
$data[11,7] = 23.01774999718594
$data[11,8] = 7.683827613233108
$data[11,9] = 9.256709345001731
$data[11,10] = 12.91686253389508
$data[11,11] = 0
$data[11,12] = 18.63159418747652
$data[11,13] = 21.93242358780391
$data[12,0] = 13.24392793181575
$data[12,1] = 10.40658925496243
$data[12,2] = 0
$data[12,3] = 16.63744090422387
$data[12,4] = 36.46432910299499
$data[12,5] = 28.99275829411849
$data[12,6] = 14.40182828336017
$data[12,7] = 23.03121612071861
$data[12,8] = 7.680948840891472
$data[12,9] = 9.219801878353067
$data[12,10] = 12.90922013902144
$data[12,11] = 0
$data[12,12] = 18.64030733919906
$data[12,13] = 21.94376608439991
$data[13,0] = 13.21141403095255
$data[13,1] = 10.40841244519884
$data[13,2] = 0
$data[13,3] = 16.63473860751798
$data[13,4] = 36.46007724129031
$data[13,5] = 28.99983392720338
$data[13,6] = 14.40615699812654
$data[13,7] = 23.03953601120284
$data[13,8] = 7.679182311417356
$data[13,9] = 9.197133436038488
$data[13,10] = 12.90457108540468
$data[13,11] = 0
$data[13,12] = 18.6456737183287
$data[13,13] = 21.95078848191232
$data[14,0] = 13.0238939282398
$data[14,1] = 10.4190489611471
$data[14,2] = 0
$data[14,3] = 16.62000001448615
$data[14,4] = 36.4380288917284
$data[14,5] = 29.04269570225657
$data[14,6] = 14.43151604421824
$data[14,7] = 23.08830985575134
$data[14,8] = 7.669010332449794
$data[14,9] = 9.066320365312571
$data[14,10] = 12.8784261615
$data[14,11] = 0
$data[14,12] = 18.67687632239011
$data[14,13] = 21.992176916394
$data[15,0] = 12.9078458236587
$data[15,1] = 10.42574268333609
$data[15,2] = 0
$data[15,3] = 16.61162813564825
$data[15,4] = 36.42657596707922
$data[15,5] = 29.07105150316727
$data[15,6] = 14.44756665440476
$data[15,7] = 23.11920894072262
$data[15,8] = 7.662727019141063
$data[15,9] = 8.985297315807911
$data[15,10] = 12.86283521475562
$data[15,11] = 0
$data[15,12] = 18.69642038620871
$data[15,13] = 22.0185894325488
$data[16,0] = 12.84074969471769
$data[16,1] = 10.42965480853554
$data[16,2] = 0
$data[16,3] = 16.60706016635255
$data[16,4] = 36.42075232264829
$data[16,5] = 29.0881171684492
$data[16,6] = 14.4569800643932
$data[16,7] = 23.13734066063496
$data[16,8] = 7.659097303536791
$data[16,9] = 8.938426225569119
$data[16,10] = 12.85403422813104
$data[16,11] = 0
$data[16,12] = 18.70780972314325
$data[16,13] = 22.03415662742027
$data[17,0] = 12.81797451456564
$data[17,1] = 10.43099006008552
$data[17,2] = 0
$data[17,3] = 16.6055561151377
$data[17,4] = 36.41891188913425
$data[17,5] = 29.09402504043712
$data[17,6] = 14.4601984702992
$data[17,7] = 23.14354147455936
$data[17,8] = 7.657865656681071
$data[17,9] = 8.922511822776386
$data[17,10] = 12.85108313476916
$data[17,11] = 0
$data[17,12] = 18.71169142071871
$data[17,13] = 22.03949187206969
$data[18,0] = 12.92023592271378
$data[18,1] = 10.4250237036974
$data[18,2] = 0
$data[18,3] = 16.61249376632579
$data[18,4] = 36.42771614454592
$data[18,5] = 29.06795468714911
$data[18,6] = 14.44583925495782
$data[18,7] = 23.11588248567255
$data[18,8] = 7.663397518480817
$data[18,9] = 8.993950530421634
$data[18,10] = 12.86447770488736
$data[18,11] = 0
$data[18,12] = 18.69432456476913
$data[18,13] = 22.01573891692792
$data[19,0] = 13.25949729353084
$data[19,1] = 10.4057183165482
$data[19,2] = 0
$data[19,3] = 16.63874990197118
$data[19,4] = 36.46640890149776
$data[19,5] = 28.989406319703
$data[19,6] = 14.39976232424046
$data[19,7] = 23.02724589306187
$data[19,8] = 7.681795073372494
$data[19,9] = 9.230655403844532
$data[19,10] = 12.91145820959881
$data[19,11] = 0
$data[19,12] = 18.63774196612106
$data[19,13] = 21.940418930841
$data[20,0] = 13.47764411492022
$data[20,1] = 10.39365120868003
$data[20,2] = 0
$data[20,3] = 16.65808646506542
$data[20,4] = 36.49845449365601
$data[20,5] = 28.94483344352833
$data[20,6] = 14.37126259562601
$data[20,7] = 22.9725169788207
$data[20,8] = 7.693676146154864
$data[20,9] = 9.382641693990296
$data[20,10] = 12.94360590898344
$data[20,11] = 0
$data[20,12] = 18.60207508703161
$data[20,13] = 21.8945399262748
$data[21,0] = 13.36156482246041
$data[21,1] = 10.40004144771311
$data[21,2] = 0
$data[21,3] = 16.64756830856316
$data[21,4] = 36.48073514104028
$data[21,5] = 28.9680026556779
$data[21,6] = 14.38632584335181
$data[21,7] = 23.00143400788064
$data[21,8] = 7.687348196681787
$data[21,9] = 9.301786942025506
$data[21,10] = 12.92631798848594
$data[21,11] = 0
$data[21,12] = 18.62099129770932
$data[21,13] = 21.91871997847829
$data[22,0] = 12.91463553250823
$data[22,1] = 10.42534855558839
$data[22,2] = 0
$data[22,3] = 16.612101650856
$data[22,4] = 36.4271983000613
$data[22,5] = 29.06935237965585
$data[22,6] = 14.44661963377861
$data[22,7] = 23.11738523150153
$data[22,8] = 7.663094439757784
$data[22,9] = 8.990039311286051
$data[22,10] = 12.8637346286361
$data[22,11] = 0
$data[22,12] = 18.69527160838586
$data[22,13] = 22.01702644428737
$data[23,0] = 12.41956171814241
$data[23,1] = 10.45502123996589
$data[23,2] = 0
$data[23,3] = 16.58331257440131
$data[23,4] = 36.39866636361211
$data[23,5] = 29.20756269891836
$data[23,6] = 14.51859778482419
$data[23,7] = 23.25619778057917
$data[23,8] = 7.645427056532669
$data[23,9] = 8.761559265798903
$data[23,10] = 12.82233350930811
$data[23,11] = 0
$data[23,12] = 18.75140193843353
$data[23,13] = 22.09494643701044

$ws.Range("B2:O25").Value = $data
